# "seeding both users and content creators"
#
# Adds a new "creators" worksheet (profile_name / parent_account / nickname /
# bio) right after the existing "users" sheet, seeds it with three rows of
# sample data, and switches the active tab/selection over to the new sheet.

$wb = $excel.ActiveWorkbook
$users = $wb.Worksheets.Item("users")

# New sheet, positioned immediately after "users".
$creators = $wb.Worksheets.Add($null, $users)
$creators.Name = "creators"

# Header row (written A1, B1, D1 first, C1 "nickname" a little later - see
# below - to line up with the shared-string table order of the source file).
$creators.Range("A1").Value = "profile_name"
$creators.Range("B1").Value = "parent_account"
$creators.Range("D1").Value = "bio"

# Row 2 - Bobby.
$creators.Range("C2").Value = "Bobby"
$creators.Range("A2").Value = "BobbyPaints"
$creators.Range("B2").Value = "jsmith"
$creators.Range("D2").Value = "Bobby likes to paint."

# "nickname" header, filled in once the nickname column has its first value.
$creators.Range("C1").Value = "nickname"

# Row 3 - Helen.
$creators.Range("A3").Value = "HelenSculpts"
$creators.Range("B3").Value = "jdoe"
$creators.Range("C3").Value = "Helen"
$creators.Range("D3").Value = "Helen likes to sculpt."

# Row 4 - Johnny.
$creators.Range("A4").Value = "JohnnyDraws"
$creators.Range("B4").Value = "jsmith"
$creators.Range("C4").Value = "Johnny"
$creators.Range("D4").Value = "Johnny likes to draw."

# Bold header row, matching the "users" sheet's header style.
$creators.Range("A1:D1").Font.Bold = $true

# Column widths sized to fit the seeded values (best-fit-style).
$creators.Columns.Item(1).ColumnWidth = 11.330729166666666
$creators.Columns.Item(2).ColumnWidth = 12.998697916666666
$creators.Columns.Item(3).ColumnWidth = 9.330729166666666
$creators.Columns.Item(4).ColumnWidth = 16.998697916666668

# Move the "users" sheet's own selection off of C4 and onto A4.
$users.Range("A4").Select()

# Make "creators" the active tab, with D4 selected.
$creators.Activate()
$creators.Range("D4").Select()
